$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.400.83"
$ws.Range("D3").Value = "'1.569.52"
$ws.Range("E3").Value = "  -4.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'291.58"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("D7").Value = "'0.3666"
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("D8").Value = "'49.42"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "'0.3386"
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("D10").Value = "'1.174"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("D11").Value = "'0.07594"
$ws.Range("E11").Value = "  -6.07%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'21.21"
$ws.Range("E13").Value = "  -4.00%  "
$ws.Range("D14").Value = "'6.065"
$ws.Range("E14").Value = "  -5.13%  "
$ws.Range("D15").Value = "'6.901"
$ws.Range("E15").Value = "  -5.98%  "
$ws.Range("D16").Value = "'0.00001140"
$ws.Range("E16").Value = "  -5.06%  "
$ws.Range("D17").Value = "'1.567.98"
$ws.Range("E17").Value = "  -4.78%  "
$ws.Range("D18").Value = "'89.11"
$ws.Range("E18").Value = "  -8.16%  "
$ws.Range("D19").Value = "'0.06754"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'6.240"
$ws.Range("E21").Value = "  -7.65%  "
$ws.Range("D22").Value = "'0.5307"
$ws.Range("E22").Value = "  -7.67%  "
$ws.Range("D23").Value = "'16.45"
$ws.Range("E23").Value = "  -5.52%  "
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").Value = "'22.406.81"
$ws.Range("E25").Value = "  -4.57%  "
$ws.Range("D26").Value = "'2.393"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("D27").Value = "'3.001"
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("D28").Value = "'19.88"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("D29").Value = "'144.89"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("D30").Value = "'4.960"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Value = "'125.33"
$ws.Range("E31").Value = "  -5.65%  "
$ws.Range("D32").Value = "'1.744.13"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").Value = "'1.043"
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("D34").Value = "'6.281"
$ws.Range("E34").Value = "  -9.67%  "
$ws.Range("D35").Value = "'1.982"
$ws.Range("E35").Value = "  -7.77%  "
$ws.Range("D36").Value = "'10.37"
$ws.Range("E36").Value = "  -9.46%  "
$ws.Range("D37").Value = "'0.02566"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").Value = "'0.08435"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'0.2305"
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("D40").Value = "'0.06542"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("D41").Value = "'5.539"
$ws.Range("E41").Value = "  -6.78%  "
$ws.Range("D42").Value = "'11.89"
$ws.Range("E42").Value = "  -8.14%  "
$ws.Range("D43").Value = "'1.255"
$ws.Range("D44").Value = "'0.6400"
$ws.Range("E44").Value = "  -7.04%  "
$ws.Range("D45").Value = "'14.43"
$ws.Range("E45").Value = "  -7.95%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'0.6026"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("D48").Value = "'3.782"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").Value = "'123.04"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "'1.216"
$ws.Range("E51").Value = "  +2.81%  "